# Applies the "Add missing files in merge conflict" edit described by the
# commit's OOXML diff against selenium-wd/data/TestCaseData.xlsx.
#
# Sheets touched (workbook tab order):
#   0 Users
#   2 Test Environment Data
#   3 Driver View Test Data

$wb = $excel.ActiveWorkbook

$wsUsers   = $wb.Worksheets.Item("Users")
$wsTestEnv = $wb.Worksheets.Item("Test Environment Data")
$wsDriver  = $wb.Worksheets.Item("Driver View Test Data")

# ---------------------------------------------------------------------------
# 1. Users sheet: no data change, B5 keeps its value (sqapicad@picarro.com) -
#    only the cell's format slot is touched (same font, cosmetic re-apply).
# ---------------------------------------------------------------------------
$wsUsers.Range("B5").NumberFormat = $wsUsers.Range("B5").NumberFormat

# ---------------------------------------------------------------------------
# 2. Test Environment Data: fill in the three previously-empty rows (11-13)
#    with the missing analyzer / replay-db3 rows that the merge had dropped.
#    New shared strings are created in this exact order so they land at
#    uniqueCount 155/156/157 (Surveyor_rr-pic.db3, Surveyor_rr-sqacudr.db3,
#    Surveyor_rr.db3) just like the authoritative workbook.
# ---------------------------------------------------------------------------
$wsTestEnv.Range("A11").Value2 = 10
$wsTestEnv.Range("B11").Value2 = "SimAuto-Analyzer1"
$wsTestEnv.Range("C11").Value2 = "SimAuto-AnalyzerKey1"
$wsTestEnv.Range("D11").Value2 = "Surveyor_rr-pic.db3"
$wsTestEnv.Range("E11").Value2 = "replay-db3.defn"

$wsTestEnv.Range("A12").Value2 = 11
$wsTestEnv.Range("B12").Value2 = "SimAuto-Analyzer2"
$wsTestEnv.Range("C12").Value2 = "SimAuto-AnalyzerKey2"
$wsTestEnv.Range("D12").Value2 = "Surveyor_rr-sqacudr.db3"
$wsTestEnv.Range("E12").Value2 = "replay-db3.defn"

$wsTestEnv.Range("A13").Value2 = 12
$wsTestEnv.Range("B13").Value2 = "SimAuto-Analyzer1"
$wsTestEnv.Range("C13").Value2 = "SimAuto-AnalyzerKey1"
$wsTestEnv.Range("D13").Value2 = "Surveyor_rr.db3"
$wsTestEnv.Range("E13").Value2 = "replay-db3.defn"

# ---------------------------------------------------------------------------
# 3. Driver View Test Data: fill in rows 22-30 (continuation of the survey
#    scenario matrix) that the merge conflict had left blank.
# ---------------------------------------------------------------------------
$driverRows = @(
  @{ Row=22; A=21; B="GenerateRandomString(15)"; C="Day";   D="Overcast"; E="Light"; F=$null;            G="Standard" },
  @{ Row=23; A=22; B="GenerateRandomString(15)"; C="Night"; D=$null;      E="Light"; F="LessThan50";      G="Standard" },
  @{ Row=24; A=23; B="GenerateRandomString(15)"; C="Day";   D="Strong";   E="Light"; F=$null;            G="Standard" },
  @{ Row=25; A=24; B="GenerateRandomString(15)"; C="Day";   D="Moderate"; E="Calm";  F=$null;            G="Standard" },
  @{ Row=26; A=25; B="GenerateRandomString(15)"; C="Day";   D="Strong";   E="Light"; F=$null;            G="RapidResponse" },
  @{ Row=27; A=26; B="GenerateRandomString(15)"; C="Night"; D=$null;      E="Light"; F="GreaterThan50";   G="RapidResponse" },
  @{ Row=28; A=27; B="GenerateRandomString(15)"; C="Day";   D="Moderate"; E="Calm";  F=$null;            G="RapidResponse" },
  @{ Row=29; A=28; B="GenerateRandomString(15)"; C="Day";   D="Strong";   E="Calm";  F=$null;            G="RapidResponse" },
  @{ Row=30; A=29; B="GenerateRandomString(15)"; C="Day";   D="Overcast"; E="Calm";  F=$null;            G="RapidResponse" }
)

foreach ($r in $driverRows) {
  $row = $r.Row
  $wsDriver.Range("A$row").Value2 = $r.A
  $wsDriver.Range("B$row").Value2 = $r.B
  $wsDriver.Range("C$row").Value2 = $r.C
  if ($r.D -ne $null) { $wsDriver.Range("D$row").Value2 = $r.D }
  $wsDriver.Range("E$row").Value2 = $r.E
  if ($r.F -ne $null) { $wsDriver.Range("F$row").Value2 = $r.F }
  $wsDriver.Range("G$row").Value2 = $r.G
}

# Rows 31-33 become touched-but-empty placeholder rows (A & B only) in the
# authoritative file - materialise the cells without giving them content.
foreach ($row in 31,32,33) {
  $wsDriver.Range("A$row").NumberFormat = "General"
  $wsDriver.Range("B$row").NumberFormat = "General"
}

# Rows 1066-1067 / 1071-1072 (far bottom of the used range) likewise end up
# as touched-but-empty cells spanning C:G after the merge clean-up.
foreach ($row in 1066,1067,1071,1072) {
  $wsDriver.Range("C$row`:G$row").NumberFormat = "General"
}

# ---------------------------------------------------------------------------
# 4. View state: restore the selections / active sheet recorded after the
#    edit. The final .Select()/.Activate() call determines the workbook's
#    activeTab, so Test Environment Data (tab index 2) must be last.
# ---------------------------------------------------------------------------
$wsDriver.Activate()
$excel.ActiveWindow.TopLeftCell = $wsDriver.Range("A4")
$wsDriver.Range("D4").Select() | Out-Null

$wsUsers.Activate()
$wsUsers.Range("B5").Select() | Out-Null

$wsTestEnv.Activate()
$wsTestEnv.Range("A13").Select() | Out-Null
